$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell R1: "tags", styled like the other bold header cells
$ws.Cells.Item(1, 18).Value = "tags"
$ws.Cells.Item(1, 18).Font.Bold = $true
$ws.Cells.Item(1, 18).HorizontalAlignment = -4108
$ws.Cells.Item(1, 18).VerticalAlignment = -4160
$ws.Cells.Item(1, 18).Borders.LineStyle = 1

# Data cells R2:R100: tags extracted per tweet
$ws.Cells.Item(2, 18).Value = "['Lula']"
$ws.Cells.Item(3, 18).Value = "['Lula']"
$ws.Cells.Item(4, 18).Value = "['Lula']"
$ws.Cells.Item(5, 18).Value = "['Lula']"
$ws.Cells.Item(6, 18).Value = "['Lula']"
$ws.Cells.Item(7, 18).Value = "['Lula']"
$ws.Cells.Item(8, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(9, 18).Value = "['Lula']"
$ws.Cells.Item(10, 18).Value = "['Lula']"
$ws.Cells.Item(11, 18).Value = "['Lula']"
$ws.Cells.Item(12, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(13, 18).Value = "['Lula']"
$ws.Cells.Item(14, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(15, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(16, 18).Value = "['Lula']"
$ws.Cells.Item(17, 18).Value = "['Lula']"
$ws.Cells.Item(18, 18).Value = "['Lula']"
$ws.Cells.Item(19, 18).Value = "['Lula']"
$ws.Cells.Item(20, 18).Value = "['Lula']"
$ws.Cells.Item(21, 18).Value = "['Lula']"
$ws.Cells.Item(22, 18).Value = "['Lula']"
$ws.Cells.Item(23, 18).Value = "['Lula']"
$ws.Cells.Item(24, 18).Value = "['Lula']"
$ws.Cells.Item(25, 18).Value = "['Lula']"
$ws.Cells.Item(26, 18).Value = "['Lula']"
$ws.Cells.Item(27, 18).Value = "['Lula']"
$ws.Cells.Item(28, 18).Value = "['Lula']"
$ws.Cells.Item(29, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(30, 18).Value = "['Lula']"
$ws.Cells.Item(31, 18).Value = "['Lula']"
$ws.Cells.Item(32, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(33, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(34, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(35, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(36, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(37, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(38, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(39, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(40, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(41, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(42, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(43, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(44, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(45, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(46, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(47, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(48, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(49, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(50, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(51, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(52, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(53, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(54, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(55, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(56, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(57, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(58, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(59, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(60, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(61, 18).Value = "[]"
$ws.Cells.Item(62, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(63, 18).Value = "[]"
$ws.Cells.Item(64, 18).Value = "[]"
$ws.Cells.Item(65, 18).Value = "[]"
$ws.Cells.Item(66, 18).Value = "[]"
$ws.Cells.Item(67, 18).Value = "[]"
$ws.Cells.Item(68, 18).Value = "['Lula']"
$ws.Cells.Item(69, 18).Value = "[]"
$ws.Cells.Item(70, 18).Value = "[]"
$ws.Cells.Item(71, 18).Value = "[]"
$ws.Cells.Item(72, 18).Value = "[]"
$ws.Cells.Item(73, 18).Value = "[]"
$ws.Cells.Item(74, 18).Value = "[]"
$ws.Cells.Item(75, 18).Value = "[]"
$ws.Cells.Item(76, 18).Value = "[]"
$ws.Cells.Item(77, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(78, 18).Value = "[]"
$ws.Cells.Item(79, 18).Value = "['Lula', 'Bolsonaro']"
$ws.Cells.Item(80, 18).Value = "[]"
$ws.Cells.Item(81, 18).Value = "[]"
$ws.Cells.Item(82, 18).Value = "[]"
$ws.Cells.Item(83, 18).Value = "[]"
$ws.Cells.Item(84, 18).Value = "[]"
$ws.Cells.Item(85, 18).Value = "[]"
$ws.Cells.Item(86, 18).Value = "[]"
$ws.Cells.Item(87, 18).Value = "[]"
$ws.Cells.Item(88, 18).Value = "[]"
$ws.Cells.Item(89, 18).Value = "[]"
$ws.Cells.Item(90, 18).Value = "[]"
$ws.Cells.Item(91, 18).Value = "[]"
$ws.Cells.Item(92, 18).Value = "[]"
$ws.Cells.Item(93, 18).Value = "[]"
$ws.Cells.Item(94, 18).Value = "['Bolsonaro']"
$ws.Cells.Item(95, 18).Value = "[]"
$ws.Cells.Item(96, 18).Value = "[]"
$ws.Cells.Item(97, 18).Value = "[]"
$ws.Cells.Item(98, 18).Value = "['Lula']"
$ws.Cells.Item(99, 18).Value = "['Lula']"
$ws.Cells.Item(100, 18).Value = "[]"

Write-Output "done"
